$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.317.75"
$ws.Range("E2").Value = "  -3.70%  "
$ws.Range("D3").Value = "3.152.38"
$ws.Range("E3").Value = "  -3.33%  "
$ws.Range("E4").Value = "  +0.21%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "607.03"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "147.39"
$r.Style = "Normal"
$ws.Range("E6").Value = "  -7.05%  "
$ws.Range("D8").Value = "3.150.57"
$ws.Range("E8").Value = "  -3.43%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.524"
$r.Style = "Normal"
$ws.Range("E9").Value = "  -4.58%  "
$ws.Range("E11").Value = "  -6.76%  "
$ws.Range("E12").Value = "  -5.93%  "
$ws.Range("E13").Value = "  -7.83%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "35.66"
$r.Style = "Normal"
$ws.Range("E14").Value = "  -9.44%  "
$ws.Range("D15").Value = "3.674.09"
$ws.Range("E15").Value = "  -3.19%  "
$ws.Range("D16").Value = "64.323.76"
$ws.Range("E16").Value = "  -3.76%  "
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "3.157.94"
$ws.Range("E18").Value = "  -2.55%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "6.93"
$r.Style = "Normal"
$ws.Range("E19").Value = "  -6.93%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "480.59"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -5.69%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "14.75"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -4.52%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "0.711"
$r.Style = "Normal"
$ws.Range("E22").Value = "  -5.59%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "13.70"
$r.Style = "Normal"
$ws.Range("E24").Value = "  -8.06%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "83.69"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -3.52%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -5.23%  "
$ws.Range("E28").Value = "  -7.15%  "
$ws.Range("E29").Value = "  -9.38%  "
$ws.Range("E30").Value = "  -1.09%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "0.114"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -18.83%  "
$ws.Range("E32").Value = "  -6.10%  "
$ws.Range("E33").Value = "  +0.02%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "26.23"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -6.86%  "
$ws.Range("E35").Value = "  -4.65%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "54.61"
$r.Style = "Normal"
$ws.Range("E36").Value = "  -2.15%  "
$ws.Range("E37").Value = "  -7.97%  "
$ws.Range("E38").Value = "  -8.55%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "453.54"
$r.Style = "Normal"
$ws.Range("E39").Value = "  -8.60%  "
$ws.Range("E40").Value = "  -13.90%  "
$ws.Range("E41").Value = "  -7.45%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.119"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -8.03%  "
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "8.44"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -4.89%  "
$ws.Range("D44").Value = "2.850.39"
$ws.Range("E44").Value = "  -4.37%  "
$ws.Range("E45").Value = "  -10.08%  "
$ws.Range("E46").Value = "  -9.89%  "
$ws.Range("E47").Value = "  -8.19%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("E49").Value = "  -7.78%  "
$ws.Range("E50").Value = "  -4.91%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "119.80"
$r.Style = "Normal"
$ws.Range("E51").Value = "  -1.60%  "
